# Generate Report for handoff
# This script updates the localization-status workbook:
#  - renames "test-content-1.md" references to a new guid-named file
#    "ab742031-5dec-45c3-b2d6-83c8f8c9d80e.md"
#  - adds a brand new tracked file
#    "eaac93a3-0636-47e9-9ac0-669f7de648ed.md" as a new row on every sheet
#  - refreshes handoff timestamps / xlf handoff filenames for the renamed file
#  - adds matching handoff info for the new file
#  - keeps the ".localization-config" row as the last row on every sheet

$wb = $excel.ActiveWorkbook

$baseRepo   = "https://github.com/OpenLocalizationTest/oltest/blob/360d8ad5898d7c111e52e0c312fa785254298175"
$handoffZh  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a1468c12224f77a8bd11c4c163ed7734339a846d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho"
$handoffDe  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/979b628de2aab031c220eb92fd4b7705db9ff2b9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho"

$file1 = "ab742031-5dec-45c3-b2d6-83c8f8c9d80e.md"
$file2 = "eaac93a3-0636-47e9-9ac0-669f7de648ed.md"

$xlf1zh = "ab742031-5dec-45c3-b2d6-83c8f8c9d80e.de557bf5b8765db9f9e69f1af90bd52d80454c9d.zh-cn.xlf"
$xlf2zh = "eaac93a3-0636-47e9-9ac0-669f7de648ed.74512ce61435bd3a1520414004d56ae47b0d0f73.zh-cn.xlf"
$xlf1de = "ab742031-5dec-45c3-b2d6-83c8f8c9d80e.de557bf5b8765db9f9e69f1af90bd52d80454c9d.de-de.xlf"
$xlf2de = "eaac93a3-0636-47e9-9ac0-669f7de648ed.74512ce61435bd3a1520414004d56ae47b0d0f73.de-de.xlf"

$readyForHandoff = "Ready for handoff"
$notToBeLocalized = "Not to be localized"
$configFile = ".localization-config"
$zeroDate = "0001-01-01 00:00:00"
$include = "Include"
$ignored = "Ignored"

$timeZh = "2016-01-08 17:34:02"
$timeDe = "2016-01-08 17:34:13"

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = 15570276
}

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Cells.Hyperlinks.Delete()

# insert a fresh row for the new file right above the config row, pushing the
# config row down from row 3 to row 4 (styles/format are inherited from the
# surrounding rows automatically)
$ws1.Rows.Item(3).Insert()

$ws1.Range("A2").Value = $file1
$ws1.Range("B2").Value = $readyForHandoff
$ws1.Range("C2").Value = $readyForHandoff

$ws1.Range("A3").Value = $file2
$ws1.Range("B3").Value = $readyForHandoff
$ws1.Range("C3").Value = $readyForHandoff

$ws1.Range("A4").Value = $configFile
$ws1.Range("B4").Value = $notToBeLocalized
$ws1.Range("C4").Value = $notToBeLocalized

$ws1.Hyperlinks.Add($ws1.Range("A2"), "$baseRepo/e2e/$file1", $null, $null, $file1) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "$baseRepo/e2e/$file2", $null, $null, $file2) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "$baseRepo/$configFile", $null, $null, $configFile) | Out-Null

Style-AsHyperlink $ws1.Range("A2")
Style-AsHyperlink $ws1.Range("A3")
Style-AsHyperlink $ws1.Range("A4")

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Cells.Hyperlinks.Delete()

$ws2.Rows.Item(3).Insert()

$ws2.Range("A2").Value = $file1
$ws2.Range("B2").Value = $readyForHandoff
$ws2.Range("C2").Value = $xlf1zh
$ws2.Range("D2").Value = $timeZh
$ws2.Range("G2").Value = $zeroDate
$ws2.Range("H2").Value = $include

$ws2.Range("A3").Value = $file2
$ws2.Range("B3").Value = $readyForHandoff
$ws2.Range("C3").Value = $xlf2zh
$ws2.Range("D3").Value = $timeZh
$ws2.Range("G3").Value = $zeroDate
$ws2.Range("H3").Value = $include

$ws2.Range("A4").Value = $configFile
$ws2.Range("B4").Value = $notToBeLocalized
$ws2.Range("D4").Value = $zeroDate
$ws2.Range("G4").Value = $zeroDate
$ws2.Range("H4").Value = $ignored

$ws2.Hyperlinks.Add($ws2.Range("A2"), "$baseRepo/e2e/$file1", $null, $null, $file1) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "$handoffZh/$xlf1zh", $null, $null, $xlf1zh) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "$baseRepo/e2e/$file2", $null, $null, $file2) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "$handoffZh/$xlf2zh", $null, $null, $xlf2zh) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "$baseRepo/$configFile", $null, $null, $configFile) | Out-Null

Style-AsHyperlink $ws2.Range("A2")
Style-AsHyperlink $ws2.Range("C2")
Style-AsHyperlink $ws2.Range("A3")
Style-AsHyperlink $ws2.Range("C3")
Style-AsHyperlink $ws2.Range("A4")

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Cells.Hyperlinks.Delete()

$ws3.Rows.Item(3).Insert()

$ws3.Range("A2").Value = $file1
$ws3.Range("B2").Value = $readyForHandoff
$ws3.Range("C2").Value = $xlf1de
$ws3.Range("D2").Value = $timeDe
$ws3.Range("G2").Value = $zeroDate
$ws3.Range("H2").Value = $include

$ws3.Range("A3").Value = $file2
$ws3.Range("B3").Value = $readyForHandoff
$ws3.Range("C3").Value = $xlf2de
$ws3.Range("D3").Value = $timeDe
$ws3.Range("G3").Value = $zeroDate
$ws3.Range("H3").Value = $include

$ws3.Range("A4").Value = $configFile
$ws3.Range("B4").Value = $notToBeLocalized
$ws3.Range("D4").Value = $zeroDate
$ws3.Range("G4").Value = $zeroDate
$ws3.Range("H4").Value = $ignored

$ws3.Hyperlinks.Add($ws3.Range("A2"), "$baseRepo/e2e/$file1", $null, $null, $file1) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "$handoffDe/$xlf1de", $null, $null, $xlf1de) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "$baseRepo/e2e/$file2", $null, $null, $file2) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "$handoffDe/$xlf2de", $null, $null, $xlf2de) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "$baseRepo/$configFile", $null, $null, $configFile) | Out-Null

Style-AsHyperlink $ws3.Range("A2")
Style-AsHyperlink $ws3.Range("C2")
Style-AsHyperlink $ws3.Range("A3")
Style-AsHyperlink $ws3.Range("C3")
Style-AsHyperlink $ws3.Range("A4")

Write-Host "Localization status workbook updated."
